# Auto-generated edit script: applies numeric updates to Sheets per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 71430664
$ws.Range("I18").Value = 1921.5834
$ws.Range("K18").Value = 1921.5834
$ws.Range("M18").Value = -1637.5834
$ws.Range("H19").Value = 743.6667
$ws.Range("I19").Value = 234
$ws.Range("J19").Value = 1253.3334
$ws.Range("K19").Value = 234
$ws.Range("L19").Value = 1253.3334
$ws.Range("M19").Value = -59
$ws.Range("N19").Value = -1603.3334
$ws.Range("H40").Value = 4763.6665
$ws.Range("J40").Value = 4358.625
$ws.Range("L40").Value = 4358.625
$ws.Range("N40").Value = -4708.625
$ws.Range("H64").Value = 125007360
$ws.Range("I64").Value = 8414
$ws.Range("K64").Value = 8414
$ws.Range("M64").Value = -8166
$ws.Range("H67").Value = 125007360
$ws.Range("I67").Value = 8414
$ws.Range("K67").Value = 8414
$ws.Range("M67").Value = -7556
$ws.Range("H112").Value = 2459.75
$ws.Range("J112").Value = 2459.75
$ws.Range("L112").Value = 7379.25
$ws.Range("N112").Value = -9595.25
$ws.Range("H116").Value = 3601.6667
$ws.Range("I116").Value = 2057.3333
$ws.Range("K116").Value = 2057.3333
$ws.Range("M116").Value = 1384.6667
$ws.Range("H131").Value = 3335077.8
$ws.Range("I131").Value = 5002217
$ws.Range("J131").Value = 799
$ws.Range("K131").Value = 15006651
$ws.Range("L131").Value = 2397
$ws.Range("M131").Value = -15001611
$ws.Range("N131").Value = -12477
$ws.Range("H132").Value = 8442.6
$ws.Range("I132").Value = 6711.0835
$ws.Range("J132").Value = 49999
$ws.Range("K132").Value = 20133.2505
$ws.Range("L132").Value = 149997
$ws.Range("M132").Value = -17603.2505
$ws.Range("N132").Value = -155057
$ws.Range("H141").Value = 6475.2
$ws.Range("I141").Value = 3861.3333
$ws.Range("K141").Value = 11583.9999
$ws.Range("M141").Value = -6403.999899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 73500
$ws.Range("J44").Value = 73500
$ws.Range("L44").Value = 73500
$ws.Range("N44").Value = -74476
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null
$ws.Range("H61").Value = 4747.154
$ws.Range("I61").Value = 2962.25
$ws.Range("K61").Value = 2962.25
$ws.Range("M61").Value = -2750.25
$ws.Range("H122").Value = 3809.366
$ws.Range("I122").Value = 3504.4856
$ws.Range("K122").Value = 10513.4568
$ws.Range("M122").Value = -8063.4568
$ws.Range("H132").Value = 2963.842
$ws.Range("I132").Value = 1831.9615
$ws.Range("J132").Value = 5416.25
$ws.Range("K132").Value = 5495.8845
$ws.Range("L132").Value = 16248.75
$ws.Range("M132").Value = -2965.8845
$ws.Range("N132").Value = -21308.75
$ws.Range("H136").Value = 4747.154
$ws.Range("I136").Value = 2962.25
$ws.Range("K136").Value = 8886.75
$ws.Range("M136").Value = -6336.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = $null
$ws.Range("H20").Value = 22731774
$ws.Range("I20").Value = 31254694
$ws.Range("K20").Value = 31254694
$ws.Range("M20").Value = -31254447
$ws.Range("H105").Value = 11820903
$ws.Range("I105").Value = 835360.7
$ws.Range("J105").Value = 25003554
$ws.Range("K105").Value = 835360.7
$ws.Range("L105").Value = 25003554
$ws.Range("M105").Value = -833613.7
$ws.Range("N105").Value = -25007048
$ws.Range("H107").Value = 1545.5
$ws.Range("I107").Value = 1461.125
$ws.Range("J107").Value = 1883
$ws.Range("K107").Value = 1461.125
$ws.Range("L107").Value = 1883
$ws.Range("M107").Value = 458.875
$ws.Range("N107").Value = -5723
$ws.Range("H134").Value = 2640.1667
$ws.Range("I134").Value = 2266.1365
$ws.Range("K134").Value = 6798.4095
$ws.Range("M134").Value = -4263.4095
$ws.Range("H140").Value = 59567
$ws.Range("J140").Value = 59567
$ws.Range("L140").Value = 59567
$ws.Range("N140").Value = -69927
$ws.Range("H141").Value = 72623.836
$ws.Range("J141").Value = 79148.60000000001
$ws.Range("L141").Value = 79148.60000000001
$ws.Range("N141").Value = -89508.60000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 845.2222
$ws.Range("I22").Value = 875.875
$ws.Range("K22").Value = 875.875
$ws.Range("M22").Value = -525.875
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H94").Value = 1770.4286
$ws.Range("J94").Value = 1761.3636
$ws.Range("L94").Value = 1761.3636
$ws.Range("N94").Value = -2663.3636
$ws.Range("H132").Value = 3877.8635
$ws.Range("I132").Value = 3841.9412
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11525.8236
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -8995.8236
$ws.Range("N132").Value = -17060
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
$ws.Range("H134").Value = 4814.647
$ws.Range("I134").Value = 4523.2666
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 13569.7998
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -11034.7998
$ws.Range("N134").Value = -26070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 9140
$ws.Range("J130").Value = 5750
$ws.Range("L130").Value = 17250
$ws.Range("N130").Value = -27290

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4996.1064
$ws.Range("I122").Value = 4853.6875
$ws.Range("J122").Value = 5299.933
$ws.Range("K122").Value = 14561.0625
$ws.Range("L122").Value = 15899.799
$ws.Range("M122").Value = -12111.0625
$ws.Range("N122").Value = -20799.799
$ws.Range("H126").Value = 3826.15
$ws.Range("I126").Value = 1376.6923
$ws.Range("K126").Value = 4130.0769
$ws.Range("M126").Value = -1660.0769
$ws.Range("H132").Value = 7477.875
$ws.Range("I132").Value = 2608.4285
$ws.Range("J132").Value = 11265.223
$ws.Range("K132").Value = 7825.2855
$ws.Range("L132").Value = 33795.669
$ws.Range("M132").Value = -5295.2855
$ws.Range("N132").Value = -38855.669

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3432.8948
$ws.Range("I7").Value = 3345.889
$ws.Range("J7").Value = 4999
$ws.Range("K7").Value = 3345.889
$ws.Range("L7").Value = 4999
$ws.Range("M7").Value = -3233.889
$ws.Range("N7").Value = -5223
$ws.Range("H40").Value = 4582.6
$ws.Range("I40").Value = 4824.8823
$ws.Range("J40").Value = 3833.7273
$ws.Range("K40").Value = 4824.8823
$ws.Range("L40").Value = 3833.7273
$ws.Range("M40").Value = -4688.8823
$ws.Range("N40").Value = -4105.7273
$ws.Range("H68").Value = 2438.25
$ws.Range("I68").Value = 2209
$ws.Range("K68").Value = 2209
$ws.Range("M68").Value = -1460
$ws.Range("H71").Value = 2438.25
$ws.Range("I71").Value = 2209
$ws.Range("K71").Value = 11045
$ws.Range("M71").Value = -7301
$ws.Range("H126").Value = 3432.8948
$ws.Range("I126").Value = 3345.889
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 10037.667
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -7567.667000000001
$ws.Range("N126").Value = -19937
$ws.Range("H136").Value = 3168.024
$ws.Range("I136").Value = 2812.9
$ws.Range("J136").Value = 4055.8333
$ws.Range("K136").Value = 8438.700000000001
$ws.Range("L136").Value = 12167.4999
$ws.Range("M136").Value = -5888.700000000001
$ws.Range("N136").Value = -17267.4999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
$ws.Range("H81").Value = 2937.25
$ws.Range("I81").Value = 2456.8572
$ws.Range("K81").Value = 4913.7144
$ws.Range("M81").Value = -3852.7144
$ws.Range("H84").Value = 2937.25
$ws.Range("I84").Value = 2456.8572
$ws.Range("K84").Value = 24568.572
$ws.Range("M84").Value = -19264.572
$ws.Range("H132").Value = 3693.913
$ws.Range("I132").Value = 3808.5264
$ws.Range("J132").Value = 3149.5
$ws.Range("K132").Value = 11425.5792
$ws.Range("L132").Value = 9448.5
$ws.Range("M132").Value = -8895.5792
$ws.Range("N132").Value = -14508.5
$ws.Range("H136").Value = 58827980
$ws.Range("I136").Value = 76924330
$ws.Range("J136").Value = 14860
$ws.Range("K136").Value = 230772990
$ws.Range("L136").Value = 44580
$ws.Range("M136").Value = -230770440
$ws.Range("N136").Value = -49680
